# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" timestamps (for zh-cn and de-de
# handoff rows) and marks those same rows' Priority as "ht" (handoff type),
# matching a refreshed localization-status report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Rows whose "Latest Handoff" timestamp / Priority get refreshed on this handoff pass.
$rows = @(7, 9, 11, 12, 13, 14)

foreach ($r in $rows) {
    # zh-cn sheet: Latest Handoff Datetime (column H) moves forward 15 seconds,
    # and Priority (column E) is now flagged "ht".
    $zhcn.Range("H$r").Value = "2016-09-02 18:22:53"
    $zhcn.Range("E$r").Value = "ht"

    # de-de sheet: same treatment, its own timestamp.
    $dede.Range("H$r").Value = "2016-09-02 18:22:58"
    $dede.Range("E$r").Value = "ht"

    # Overview sheet mirrors the "Latest HO Xliff Generate Date" (column G)
    # with the de-de timestamp value already used there.
    $overview.Range("G$r").Value = "2016-09-02 18:22:58"
}
